$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 280; this shifts the existing rows 280-312
# down to 281-313, preserving all of their data/formatting.
$ws.Rows.Item(280).Insert()

# Populate the newly inserted (blank) row 280 with the new weekly record.
$ws.Cells.Item(280, 1).Value = 10
$ws.Cells.Item(280, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(280, 3).Value = "La Araucanía"
$ws.Cells.Item(280, 4).Value = 44918
$ws.Cells.Item(280, 5).Value = 9
$ws.Cells.Item(280, 6).Value = 100112039
$ws.Cells.Item(280, 7).Value = "Ciboulette"
$ws.Cells.Item(280, 8).Value = "Sin especificar"
$ws.Cells.Item(280, 9).Value = "Primera"
$ws.Cells.Item(280, 10).Value = 65
$ws.Cells.Item(280, 11).Value = 5000
$ws.Cells.Item(280, 12).Value = 5000
$ws.Cells.Item(280, 13).Value = 5000
$ws.Cells.Item(280, 14).Value = "$/docena de atados"
$ws.Cells.Item(280, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(280, 16).Value = 1667
$ws.Cells.Item(280, 17).Value = 3
$ws.Cells.Item(280, 18).Value = "Hortaliza"
